# Fix lỗi chính tả
# Rewrites the task-description cells so each bullet point starts on its own
# line with a leading "- " marker, and fixes a couple of typos
# ("cảu" -> "của") / spacing issues along the way.
#
# Cells are set through .Formula with a leading apostrophe ('...") instead of
# .Value so that Excel records them as explicit text with a quote/text
# prefix (matching the authored workbook, where every one of these strings
# now begins with "-").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$ws.Range("H4").Formula  = "'- Đọc tài liệu .Net và Oracle. $nl- Làm các phần bài tập cuối chương: chương 1 - 6  của Oracle."

$ws.Range("D11").Formula = "'- Đọc tài liệu và làm các phần bài tập cuối chương 7-9 của Oracle. $nl- Nộp bài tập chương 1-8"
$ws.Range("E11").Formula = "'- Đọc tài liệu và làm các bài tập cuối chương 1-4 của ASP.Net. $nl- Nộp bài tập chương 9 của Oracle, chương 1-3 của ASP.Net"
$ws.Range("F11").Formula = "'- Đọc tài liệu chương 5-6 của ASP.Net. $nl- Nộp bài tập chương 5 của ASP.Net. "
$ws.Range("H11").Formula = "'- Đọc tài liệu và làm bài tập chương 7 của ASP.Net. $nl- Tham gia tổ chức 20/10 tại công ty"

$ws.Range("D18").Formula = "'- Làm bài Exam 1, Exam 2 của Oracle"
$ws.Range("E18").Formula = "'- Làm bài Exam 2, Exa m3 của Oracle"
$ws.Range("F18").Formula = "'- Làm Exam 1 của Asp.Net "
$ws.Range("H18").Formula = "'- Làm tiếp Exam 1 và làm lại bài tập của Asp.Net"
